# "Fruta / hortaliza, semanal" -- weekly price update.
# A new week's record is inserted as row 3 (pushing the existing
# rows 3-5 down to rows 4-6); the new row carries the latest
# Pepino dulce ("Primera") price quote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 3, shifting the old
# rows 3, 4 and 5 down to 4, 5 and 6 respectively.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's data.
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 44635
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 100112043
$ws.Cells.Item(3, 7).Value = "Pepino dulce"
$ws.Cells.Item(3, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 15000
$ws.Cells.Item(3, 12).Value = 16000
$ws.Cells.Item(3, 13).Value = 15500
$ws.Cells.Item(3, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(3, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 16).Value = 861
$ws.Cells.Item(3, 17).Value = 18
$ws.Cells.Item(3, 18).Value = "Hortaliza"
